# Iteration 2 Plan 2.1 - apply the "Modifications to Project Plan, Gantt
# Chart and Short Use case description" edit:
#
#   1. The TITLE field in the subtitle ("Iteration 1 Plan") is unpacked
#      from a <w:fldSimple> into explicit complex-field codes
#      (fldChar begin/separate/end + instrText) and its cached result
#      text is bumped from "Iteration 1 Plan" to "Iteration 2 Plan".
#   2. A (now current) "_GoBack" bookmark is left behind at the point of
#      that last edit - the empty Title-styled paragraph right under the
#      subtitle.
#   3/4. Because a bookmark got renumbered, the pre-existing "OLE_LINK1"
#      bookmark (wrapping the milestones table) shifts from id 0 to id 1.
#   5. The stale "_GoBack" bookmark that used to sit after "In Progress"
#      in the Gantt table is gone (Word only ever keeps the most recent
#      "_GoBack").

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Subtitle field: fldSimple -> complex field, "1" -> "2"
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(2)
$titleRange = $titlePara.Range

$fieldXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="18CBCBEA" w14:textId="41F5432A" w:rsidR="003C7438" w:rsidRPr="00350E34" w:rsidRDefault="005A4753" w:rsidP="00350E34"><w:pPr><w:pStyle w:val="Heading2"/><w:jc w:val="center"/></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> TITLE  \* MERGEFORMAT </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:t xml:space="preserve">Iteration </w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Plan</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$titleRange.InsertXML($fieldXml)

# ---------------------------------------------------------------------
# 2. Drop a fresh "_GoBack" bookmark at the empty Title-styled paragraph
#    that follows the subtitle (paragraph 3). Adding it here causes the
#    pre-existing "OLE_LINK1" bookmark to be renumbered out from under
#    id 0 automatically, matching points 3/4 below.
# ---------------------------------------------------------------------
$afterTitlePara = $d.Paragraphs.Item(3)
$d.Bookmarks.Add("_GoBack", $afterTitlePara.Range)

# ---------------------------------------------------------------------
# 5. Remove the stale "_GoBack" bookmark that used to live after the
#    "In Progress" text in the Gantt chart table (the row that starts
#    with "Produce Iteration Assessment").
# ---------------------------------------------------------------------
$scan = $d.Content
$scan.Start = 0
$matchIndex = 0
$ganttRowRange = $null
while ($scan.Find.Execute("Produce Iteration Assessment", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $matchIndex = $matchIndex + 1
    if ($matchIndex -eq 2) {
        $ganttRowRange = $scan.Duplicate()
        break
    }
    $scan.Collapse(0)
}

$probe = $d.Range($ganttRowRange.End, $ganttRowRange.End + 300)
$probe.Find.Execute("In Progress", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$inProgressCell = $d.Range($probe.Start, $probe.End + 1)

$cellXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5F6A4C8A" w14:textId="53B27D0D" w:rsidR="00695A19" w:rsidRPr="006458A2" w:rsidRDefault="00D85C49" w:rsidP="006458A2"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:cs="Arial"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:lang w:eastAsia="ja-JP"/></w:rPr><w:t>In Progress</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$inProgressCell.InsertXML($cellXml)

Write-Output "Applied Iteration 2 Plan edits."
